$wb = $excel.ActiveWorkbook
$accounts = $wb.Worksheets.Item("Accounts")

# Duplicate "Accounts" and move the copy in front of it, so the new sheet
# inherits the workbook's original sheetFormatPr/default styles instead of
# the engine's blank-worksheet defaults.
$accounts.Copy($accounts)
$new = $wb.Worksheets.Item("Accounts (2)")
$new.Name = "Marketinglist"

# Clear the copied content/formatting so we can lay out the Marketinglist
# sheet from scratch.
$new.Cells.Clear()

# Fill column A (rows 2-7) first, then B/C columns row by row, then the
# header row last, to reproduce the original shared-string insertion order.
$new.Range("A2").Value = "My List 1"
$new.Range("A3").Value = "My List 2"
$new.Range("A4").Value = "My List 1"
$new.Range("A5").Value = "My List 1"
$new.Range("A6").Value = "My List 3"
$new.Range("A7").Value = "My List 3"

$new.Range("B2").Value = "Simon"
$new.Range("C2").Value = "Meyer"
$new.Range("B3").Value = "Peter"
$new.Range("C3").Value = "Chan"
$new.Range("B4").Value = "Peter"
$new.Range("C4").Value = "Chan"
$new.Range("B5").Value = "James"
$new.Range("C5").Value = "Bond"
$new.Range("B6").Value = "Simon"
$new.Range("C6").Value = "Meyer"
$new.Range("B7").Value = "James"
$new.Range("C7").Value = "Bond"

$new.Range("A1").Value = "My Marketinglists"
$new.Range("B1").Value = "Firstname"
$new.Range("C1").Value = "Lastname"

$new.Range("A1:C1").Font.Bold = $true

$new.Columns.Item(1).ColumnWidth = 15.333333333333334
$new.Columns.Item(2).ColumnWidth = 11.833333333333334

$null = $new.Range("C4").Select()

Write-Host "done"
